$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Correct the year range for "Earthquake / volcanic eruptions / major geophysical disasters"
# (row 11, column C) from "2006-2021" to "2007-2021"
$ws.Range("C11").Value = "2007-2021"

# Reflect the last clicked cell selection
$ws.Range("C12").Select()
